$d = $word.ActiveDocument

# The document ends with:
#   "Requisitos"
#   "LOM3015: Termodinamica de Materiais (Requisito fraco)"
#   <empty paragraph>
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   "(c) 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
#   <empty paragraph>
#   <empty paragraph with pageBreakBefore>
#
# This commit removes the stray "Ver no Jupiter..." / copyright footer block
# (which was scraped along with the blank paragraph immediately preceding it),
# leaving the final blank paragraph and the page-break paragraph after it intact.

# Locate the "Ver no Jupiter..." paragraph and expand to the whole paragraph
# (so the selection includes its trailing paragraph mark).
$startRange = $d.Content
$startRange.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx") | Out-Null
$startRange.Expand(4) | Out-Null
$start = $startRange.Start

# Back up one character to also remove the blank paragraph immediately
# preceding the footer text.
$start = $start - 1

# Locate the copyright paragraph and expand it so the end includes its
# trailing paragraph mark.
$endRange = $d.Content
$endRange.Find.Execute("Original theme under Creative Commons Attribution") | Out-Null
$endRange.Expand(4) | Out-Null
$end = $endRange.End

$r = $d.Range($start, $end)
$r.Delete()
